$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Match the formatting used by other populated "even" data rows (e.g. row 24)
# before filling values, same as typing into a row already formatted like its neighbours.
$ws.Range("B24:G24").Copy()
$ws.Range("B26:G26").PasteSpecial(-4122)  # xlPasteFormats

# Fill in row 26 (24th data entry): 支出 生活费(12/21-12/31), 300, on 2017-12-22
$ws.Range("B26").Value = 24
$ws.Range("C26").Value = "支出"
$ws.Range("D26").Value = 300
$ws.Range("E26").Value = Get-Date -Year 2017 -Month 12 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("F26").Value = "生活费"
$ws.Range("G26").Value = "生活费(12/21-12/31)"

# Update the active selection as recorded in the saved view
$ws.Range("L14").Select()
